# Generate Report for Handback
# Replaces the stale e2e test artifact UUID file references with the new
# run's UUIDs, and refreshes the associated handoff/handback timestamps
# across the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$oldUuid1 = "a3673701-bd42-44b4-81d0-d3c2f37199d9"
$oldUuid2 = "ec16c88b-337a-4ed3-a246-c585ea35404a"
$newUuid1 = "7945b2a5-f081-4612-9729-8eb5f78851c0"
$newUuid2 = "ffff04920d58-557b-4ec2-85f7-a6c77be4f265"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A2").Value = "$newUuid1.md"
$ws1.Range("B2").Value = "e2e\$newUuid1.md"
$ws1.Range("G2").Value = "2016-08-27 09:01:53"

$ws1.Range("A3").Value = "$newUuid2.md"
$ws1.Range("B3").Value = "e2e\$newUuid2.md"
$ws1.Range("G3").Value = "2016-08-27 09:01:53"

# Hyperlinks keep pointing at their original targets (the .rels entries are
# untouched) - only the visible display text reflects the new file names.
$ws1Link2Address = $ws1.Hyperlinks.Item(1).Address
$ws1Link3Address = $ws1.Hyperlinks.Item(2).Address
$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("B2"), $ws1Link2Address, "", "", "e2e\$newUuid1.md")
$ws1.Hyperlinks.Add($ws1.Range("B3"), $ws1Link3Address, "", "", "e2e\$newUuid2.md")

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2Link1Address = $ws2.Hyperlinks.Item(1).Address
$ws2Link2Address = $ws2.Hyperlinks.Item(2).Address
$ws2Link3Address = $ws2.Hyperlinks.Item(3).Address
$ws2Link4Address = $ws2.Hyperlinks.Item(4).Address

$ws2.Range("A2").Value = "$newUuid1.md"
$ws2.Range("I2").Value = "$newUuid1.md"
$ws2.Range("G2").Value = "$newUuid1.62e8cad1dcd6bf6aa0348605d9897ce8587f75a2.zh-cn.xlf"
$ws2.Range("H2").Value = "2016-08-27 09:01:49"
$ws2.Range("J2").Value = "$newUuid1.62e8cad1dcd6bf6aa0348605d9897ce8587f75a2.zh-cn.xlf"
$ws2.Range("K2").Value = "2016-08-27 09:02:15"

$ws2.Range("A3").Value = "$newUuid2.md"
$ws2.Range("I3").Value = "$newUuid2.md"
$ws2.Range("G3").Value = "$newUuid1.62e8cad1dcd6bf6aa0348605d9897ce8587f75a2.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-08-27 09:01:49"
$ws2.Range("J3").Value = "$newUuid1.62e8cad1dcd6bf6aa0348605d9897ce8587f75a2.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-08-27 09:02:15"

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), $ws2Link1Address, "", "", "$newUuid1.md")
$ws2.Hyperlinks.Add($ws2.Range("I2"), $ws2Link2Address, "", "", "$newUuid1.md")
$ws2.Hyperlinks.Add($ws2.Range("A3"), $ws2Link3Address, "", "", "$newUuid2.md")
$ws2.Hyperlinks.Add($ws2.Range("I3"), $ws2Link4Address, "", "", "$newUuid2.md")

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3Link1Address = $ws3.Hyperlinks.Item(1).Address
$ws3Link2Address = $ws3.Hyperlinks.Item(2).Address
$ws3Link3Address = $ws3.Hyperlinks.Item(3).Address
$ws3Link4Address = $ws3.Hyperlinks.Item(4).Address

$ws3.Range("A2").Value = "$newUuid1.md"
$ws3.Range("I2").Value = "$newUuid1.md"
$ws3.Range("G2").Value = "$newUuid1.62e8cad1dcd6bf6aa0348605d9897ce8587f75a2.de-de.xlf"
$ws3.Range("H2").Value = "2016-08-27 09:01:53"
$ws3.Range("J2").Value = "$newUuid1.62e8cad1dcd6bf6aa0348605d9897ce8587f75a2.de-de.xlf"
$ws3.Range("K2").Value = "2016-08-27 09:02:21"

$ws3.Range("A3").Value = "$newUuid2.md"
$ws3.Range("I3").Value = "$newUuid2.md"
$ws3.Range("G3").Value = "$newUuid1.62e8cad1dcd6bf6aa0348605d9897ce8587f75a2.de-de.xlf"
$ws3.Range("H3").Value = "2016-08-27 09:01:53"
$ws3.Range("J3").Value = "$newUuid1.62e8cad1dcd6bf6aa0348605d9897ce8587f75a2.de-de.xlf"
$ws3.Range("K3").Value = "2016-08-27 09:02:21"

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), $ws3Link1Address, "", "", "$newUuid1.md")
$ws3.Hyperlinks.Add($ws3.Range("I2"), $ws3Link2Address, "", "", "$newUuid1.md")
$ws3.Hyperlinks.Add($ws3.Range("A3"), $ws3Link3Address, "", "", "$newUuid2.md")
$ws3.Hyperlinks.Add($ws3.Range("I3"), $ws3Link4Address, "", "", "$newUuid2.md")
